$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '54.861.77'
$ws.Cells.Item(2, 5).Value = '  -3.47%  '

$ws.Cells.Item(3, 4).Value = '2.344.51'
$ws.Cells.Item(3, 5).Value = '  -5.70%  '

$ws.Cells.Item(4, 5).Value = '  +0.06%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '470.64'
$ws.Cells.Item(5, 5).Value = '  -4.09%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '142.81'
$ws.Cells.Item(6, 5).Value = '  -3.29%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.999'
$ws.Cells.Item(7, 5).Value = '  +0.24%  '

$ws.Cells.Item(8, 5).Value = '  -2.88%  '

$ws.Cells.Item(9, 4).Value = '2.347.15'
$ws.Cells.Item(9, 5).Value = '  -6.16%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.0957'
$ws.Cells.Item(10, 5).Value = '  -2.59%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '5.41'
$ws.Cells.Item(11, 5).Value = '  -6.68%  '

$ws.Cells.Item(12, 5).Value = '  -4.62%  '

$ws.Cells.Item(13, 5).Value = '  +0.62%  '

$ws.Cells.Item(14, 4).Value = '2.756.60'
$ws.Cells.Item(14, 5).Value = '  -5.48%  '

$ws.Cells.Item(15, 4).Value = '54.888.79'
$ws.Cells.Item(15, 5).Value = '  -3.07%  '

$ws.Cells.Item(16, 5).Value = '  -6.32%  '

$ws.Cells.Item(17, 5).Value = '  -5.29%  '

$ws.Cells.Item(18, 4).Value = '2.351.44'
$ws.Cells.Item(18, 5).Value = '  -5.88%  '

$ws.Cells.Item(19, 5).Value = '  -1.08%  '

$ws.Cells.Item(20, 5).Value = '  -2.89%  '

$ws.Cells.Item(21, 5).Value = '  -5.81%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '0.999'
$ws.Cells.Item(22, 5).Value = '  +0.06%  '

$ws.Cells.Item(23, 5).Value = '  -4.27%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '55.80'
$ws.Cells.Item(24, 5).Value = '  -5.32%  '

$ws.Cells.Item(25, 5).Value = '  +0.25%  '

$ws.Cells.Item(26, 5).Value = '  -5.34%  '

$ws.Cells.Item(27, 5).Value = '  -5.75%  '

$ws.Cells.Item(28, 4).Value = '2.451.02'
$ws.Cells.Item(28, 5).Value = '  -5.44%  '

$ws.Cells.Item(29, 5).Value = '  -7.02%  '

$ws.Cells.Item(30, 5).Value = '  -0.05%  '

$ws.Cells.Item(31, 5).Value = '  -5.83%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '146.91'
$ws.Cells.Item(32, 5).Value = '  -1.59%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '17.93'
$ws.Cells.Item(33, 5).Value = '  -1.64%  '

$ws.Cells.Item(34, 5).Value = '  -3.64%  '

$ws.Cells.Item(35, 5).Value = '  -4.37%  '

$ws.Cells.Item(36, 5).Value = '  -5.82%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '3.54'
$ws.Cells.Item(37, 5).Value = '  -5.30%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.817'
$ws.Cells.Item(38, 5).Value = '  -5.21%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '33.39'

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.999'
$ws.Cells.Item(40, 5).Value = '  +0.26%  '

$ws.Cells.Item(41, 5).Value = '  -0.11%  '

$ws.Cells.Item(42, 5).Value = '  -5.03%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.0942'
$ws.Cells.Item(43, 5).Value = '  +2.13%  '

$ws.Cells.Item(44, 5).Value = '  -5.97%  '

$ws.Cells.Item(45, 5).Value = '  -6.90%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '10.15'
$ws.Cells.Item(46, 5).Value = '  -0.32%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '251.38'
$ws.Cells.Item(47, 5).Value = '  -3.11%  '

$ws.Cells.Item(48, 5).Value = '  -4.46%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '4.38'
$ws.Cells.Item(49, 5).Value = '  -8.42%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '16.68'
$ws.Cells.Item(50, 5).Value = '  -5.45%  '

$ws.Cells.Item(51, 4).Value = '1.773.47'
$ws.Cells.Item(51, 5).Value = '  -6.42%  '
